$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value updates for column A (z_sim), per diff
$updates = @(
    @{Row=4; Value=2},
    @{Row=9; Value=2},
    @{Row=10; Value=1},
    @{Row=12; Value=1},
    @{Row=14; Value=1},
    @{Row=15; Value=2},
    @{Row=17; Value=2},
    @{Row=23; Value=1},
    @{Row=24; Value=1},
    @{Row=25; Value=1},
    @{Row=28; Value=1},
    @{Row=29; Value=1},
    @{Row=30; Value=2},
    @{Row=31; Value=2},
    @{Row=32; Value=2},
    @{Row=36; Value=2},
    @{Row=39; Value=2},
    @{Row=41; Value=1},
    @{Row=42; Value=1},
    @{Row=43; Value=2},
    @{Row=44; Value=1},
    @{Row=47; Value=1},
    @{Row=50; Value=1},
    @{Row=58; Value=2},
    @{Row=60; Value=1},
    @{Row=63; Value=2},
    @{Row=69; Value=1},
    @{Row=70; Value=2},
    @{Row=74; Value=2},
    @{Row=75; Value=2},
    @{Row=78; Value=1},
    @{Row=81; Value=2},
    @{Row=87; Value=2},
    @{Row=88; Value=2},
    @{Row=89; Value=1},
    @{Row=94; Value=2},
    @{Row=97; Value=1},
    @{Row=98; Value=2},
    @{Row=103; Value=2},
    @{Row=108; Value=1},
    @{Row=111; Value=1},
    @{Row=119; Value=1},
    @{Row=120; Value=2},
    @{Row=123; Value=1},
    @{Row=129; Value=2},
    @{Row=131; Value=2},
    @{Row=134; Value=2},
    @{Row=139; Value=2},
    @{Row=140; Value=2},
    @{Row=142; Value=1},
    @{Row=146; Value=2},
    @{Row=147; Value=1},
    @{Row=149; Value=2},
    @{Row=154; Value=1},
    @{Row=156; Value=1},
    @{Row=157; Value=1},
    @{Row=158; Value=1},
    @{Row=159; Value=1},
    @{Row=166; Value=1},
    @{Row=168; Value=2},
    @{Row=170; Value=1},
    @{Row=172; Value=1},
    @{Row=175; Value=2},
    @{Row=176; Value=2},
    @{Row=178; Value=1},
    @{Row=179; Value=2},
    @{Row=181; Value=1},
    @{Row=183; Value=1},
    @{Row=185; Value=1},
    @{Row=187; Value=1},
    @{Row=188; Value=1},
    @{Row=190; Value=2},
    @{Row=191; Value=1},
    @{Row=196; Value=1},
    @{Row=197; Value=2},
    @{Row=198; Value=1},
    @{Row=199; Value=2}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Value
}

